$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OptimizationParameters EQ2")

# The old per-run-parameter table (B3:M6, both Pressurizing + Depressurizing
# fit columns) is replaced by a single merged "kinks and P/DP" summary cell
# per row, holding the saved coefficient matrix as text.
$ws.Range("B3:M6").Clear()

$ws.Range("I3").Value = " {[     0.6187 -119.0889 5.6971 -40.2975]}"
$ws.Range("I4").Value = "    {[87.8422 2.6166e+03 124.6828 -108.6391]}"
$ws.Range("I5").Value = "    {[   36.0653 532.8312 95.4642 -155.7685]}"
$ws.Range("I6").Value = "    {[   14.3226 139.6986 62.3122 -107.7226]}"

$ws.Range("I3:L3").HorizontalAlignment = -4108
$ws.Range("I4:L4").HorizontalAlignment = -4108
$ws.Range("I5:L5").HorizontalAlignment = -4108
$ws.Range("I6:L6").HorizontalAlignment = -4108

$ws.Range("I3:L3").Merge()
$ws.Range("I4:L4").Merge()
$ws.Range("I5:L5").Merge()
$ws.Range("I6:L6").Merge()

$ws.Range("I10").Select() | Out-Null
